# Update date for reposive (Responsive) of news section
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22 ("News" task): set Actual Finish date (C22) and progress (D22) ---
# C22 was empty with the "general" style (s=7); it becomes a date cell that
# matches the formatting already used by the sibling date cell B22 (s=24).
$ws.Range("B22").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null          # xlPasteFormats
$ws.Range("C22").Value = 44481                             # 10/12/2021
$ws.Range("D22").Value = 1

# --- Row 25 ("Responsive" task): set Start/Finish dates (B25, C25) and progress (D25) ---
# B25/C25 were empty with the "general" style (s=7); they become date cells
# that match the formatting already used on row 23/24 (s=25).
$ws.Range("B23").Copy() | Out-Null
$ws.Range("B25").PasteSpecial(-4122) | Out-Null           # xlPasteFormats
$ws.Range("B25").Value = 44481

$ws.Range("C23").Copy() | Out-Null
$ws.Range("C25").PasteSpecial(-4122) | Out-Null           # xlPasteFormats
$ws.Range("C25").Value = 44481

$ws.Range("D25").Value = 1

# --- Update the sheet's visible selection / scroll position ---
$ws.Range("D23").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
